$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Budgets")

# Update the JUN (column H) budget figures for the two data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 2

# Move/save the active selection to H4, matching the saved sheet view state.
$ws.Range("H4").Select()
